# Adding Pojo classes for Customer and Customer address Refactoring
#
# The original workbook has two sheets: "TestCases" (a trivial RunMode
# switch sheet) and "TestData" (the actual test-data table). This edit
# removes the "TestCases" sheet entirely (folding its functionality away)
# and extends "TestData" with additional address[...] columns used by the
# new Pojo-based test (createCustomerWithValidTokenUsingPojo), while
# renaming the top header of the first block accordingly.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Drop the "TestCases" sheet so "TestData" becomes the only/first sheet
# ---------------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Delete()

$ws = $wb.Worksheets.Item("TestData")
$ws.Activate()

# ---------------------------------------------------------------------
# 2. Relabel the first block's title row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "createCustomerWithValidTokenUsingPojo"

# ---------------------------------------------------------------------
# 3. Shift "preferred_locales[0]" column (old column F, rows 2-4) out to
#    column K, making room for the new address[...] columns F..J
# ---------------------------------------------------------------------
$oldF2 = $ws.Range("F2").Value()
$oldF3 = $ws.Range("F3").Value()
$oldF4 = $ws.Range("F4").Value()

$ws.Range("F2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K2").Value = $oldF2

$ws.Range("F3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Value = $oldF3

$ws.Range("F4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").Value = $oldF4

# ---------------------------------------------------------------------
# 4. Populate the new address[...] header cells F2:J2 (copy formatting
#    from the neighbouring "address[city]" header at E2)
# ---------------------------------------------------------------------
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:J2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("F2").Value = "address[country]"
$ws.Range("G2").Value = "address[line2]"
$ws.Range("H2").Value = "address[line1]"
$ws.Range("I2").Value = "address[postal_code]"
$ws.Range("J2").Value = "address[state]"

# ---------------------------------------------------------------------
# 5. Populate the new address[...] data cells for both data rows (copy
#    formatting from the neighbouring "Kyiv" cell in each row: E3 / E4)
# ---------------------------------------------------------------------
$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3:J3").PasteSpecial(-4122) | Out-Null

$ws.Range("F3").Value = "Ukraine"
$ws.Range("G3").Value = "line2"
$ws.Range("H3").Value = "line1"
$ws.Range("I3").Value = 7400
$ws.Range("J3").Value = "Kyiv obl."

$ws.Range("E4").Copy() | Out-Null
$ws.Range("F4:J4").PasteSpecial(-4122) | Out-Null

$ws.Range("F4").Value = "Ukraine"
$ws.Range("G4").Value = "line2"
$ws.Range("H4").Value = "line1"
$ws.Range("I4").Value = 7400
$ws.Range("J4").Value = "Kyiv obl."

# ---------------------------------------------------------------------
# 6. Column widths: widen column A, size the new address[...] columns
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.59              # -> raw width ~39.43
$ws.Columns.Item(5).ColumnWidth = 14.02               # -> raw width ~14.86
$ws.Columns.Item(7).ColumnWidth = 13.88               # -> raw width ~14.71
$ws.Columns.Item(8).ColumnWidth = 13.74               # -> raw width ~14.57
$ws.Columns.Item(9).ColumnWidth = 18.88               # -> raw width ~19.71
$ws.Columns.Item(10).ColumnWidth = 13.88              # -> raw width ~14.71
$ws.Columns.Item(11).ColumnWidth = 18.88              # -> raw width ~19.71

# ---------------------------------------------------------------------
# 7. Selection / active cell matches the authored file
# ---------------------------------------------------------------------
$ws.Range("B2").Select()

$wb.Save()
